# Update "想去人数" (F column) values for the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Row -> new F-column value for each affected sheet.
$updates = @{
    2  = 382
    3  = 10888
    5  = 988
    6  = 206
    7  = 1352
    8  = 8348
    9  = 48
    12 = 228
    14 = 3341
    18 = 842
    19 = 136
    20 = 1080
    21 = 290
    22 = 135
    23 = 1870
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
